$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect before writing, matching the original password.
$ws.Unprotect("D382")

# Update the confidential disclaimer date: 2021-04-21 -> 2021-04-22
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-22 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values per row
$ws.Range("D2").Value = 0.02087960388890529
$ws.Range("E2").Value = -0.01312689330191863
$ws.Range("D3").Value = 0.0187003824530608
$ws.Range("E3").Value = -0.01333591032083503
$ws.Range("D4").Value = 0.02007622556357461
$ws.Range("E4").Value = -0.02592414786365815
$ws.Range("D5").Value = 0.02014288929864329
$ws.Range("E5").Value = -0.02008652657601973
$ws.Range("D6").Value = 0.01955094745833466
$ws.Range("E6").Value = -0.008873551885629882
$ws.Range("D7").Value = 0.01984169772453781
$ws.Range("E7").Value = -0.02539062499999989
$ws.Range("D8").Value = 0.0197910975400881
$ws.Range("E8").Value = 0.004159733777038266
$ws.Range("D9").Value = 0.02023846742482609
$ws.Range("E9").Value = -0.01293753472497816
$ws.Range("D10").Value = 0.01892928804938096
$ws.Range("E10").Value = -0.007892057026476662
$ws.Range("D11").Value = 0.01968588128353391
$ws.Range("E11").Value = -0.0142288861689106
$ws.Range("D12").Value = 0.01925849043987824
$ws.Range("E12").Value = -0.0127774041694686
$ws.Range("D13").Value = 0.02102036075120391
$ws.Range("E13").Value = -0.001623903864891196
$ws.Range("D14").Value = 0.02010513995468874
$ws.Range("E14").Value = -0.00251677852348986
$ws.Range("D15").Value = 0.01909233308816339
$ws.Range("E15").Value = -0.002313743637204935
$ws.Range("D16").Value = 0.01767713427863313
$ws.Range("E16").Value = -0.003839338452451435
$ws.Range("D17").Value = 0.01801226010342115
$ws.Range("E17").Value = -0.01532801961986507
$ws.Range("D18").Value = 0.01645028060011015
$ws.Range("E18").Value = -0.003368893879842738
$ws.Range("D19").Value = 0.01490898291822107
$ws.Range("E19").Value = -0.01822222222222236
$ws.Range("D20").Value = 0.02244078021468534
$ws.Range("E20").Value = -0.004473872584108807
$ws.Range("D21").Value = 0.02164965034670162
$ws.Range("E21").Value = -0.006010016694490861
$ws.Range("D22").Value = 0.02099807257472011
$ws.Range("E22").Value = 0.002553191489361506
$ws.Range("D23").Value = 0.02010453757154053
$ws.Range("E23").Value = 0.008838951310861543
$ws.Range("D24").Value = 0.01882246543776488
$ws.Range("E24").Value = -0.009387668017921902
$ws.Range("D25").Value = 0.01880720506467687
$ws.Range("E25").Value = 0.002391527160915619
$ws.Range("D26").Value = 0.02009911612320663
$ws.Range("E26").Value = -0.01126895642270587
$ws.Range("D27").Value = 0.01814980425559599
$ws.Range("E27").Value = 0.02849872773536899
$ws.Range("D28").Value = 0.02015694490543488
$ws.Range("E28").Value = -0.0147928994082841
$ws.Range("D29").Value = 0.01748718279256394
$ws.Range("E29").Value = -0.03341370995521886
$ws.Range("D30").Value = 0.01300464899234351
$ws.Range("E30").Value = -0.00200722601364911
$ws.Range("D31").Value = 0.009635118455634179
$ws.Range("E31").Value = 0.001458789204959787
$ws.Range("D32").Value = 0.01749541536225616
$ws.Range("E32").Value = 0.005107252298263587
$ws.Range("D33").Value = 0.01990314080565533
$ws.Range("E33").Value = -0.005558806319485021
$ws.Range("D34").Value = 0.02020493476290902
$ws.Range("E34").Value = 0.01714285714285713
$ws.Range("D35").Value = 0.01796427024594701
$ws.Range("E35").Value = -0.005789909015715411
$ws.Range("D36").Value = 0.01942083269832109
$ws.Range("E36").Value = 0.03333333333333344
$ws.Range("D37").Value = 0.01798234174039334
$ws.Range("E37").Value = 0.002277904328018332
$ws.Range("D38").Value = 0.02109344990652018
$ws.Range("E38").Value = -0.005759162303665133
$ws.Range("D39").Value = 0.02302870616733923
$ws.Range("E39").Value = -0.01715959821428581
$ws.Range("D40").Value = 0.0190907267331015
$ws.Range("E40").Value = 0.007194244604316502
$ws.Range("D41").Value = 0.02139022400420544
$ws.Range("E41").Value = -0.01216581556022833
$ws.Range("D42").Value = 0.01974310768261396
$ws.Range("E42").Value = -0.005797101449275255
$ws.Range("D43").Value = 0.02022983326636841
$ws.Range("E43").Value = -0.01832276250880904
$ws.Range("D44").Value = 0.01941059218480151
$ws.Range("E44").Value = -0.003806804663335828
$ws.Range("D45").Value = 0.01869897689238164
$ws.Range("E45").Value = 0.001879194630872494
$ws.Range("D46").Value = 0.01926461466855172
$ws.Range("E46").Value = 0.00160513643659721
$ws.Range("D47").Value = 0.01908259456060065
$ws.Range("E47").Value = 0.03524472433985193
$ws.Range("D48").Value = 0.0186015916167542
$ws.Range("E48").Value = -0.006476683937823857
$ws.Range("D49").Value = 0.01704141926288784
$ws.Range("E49").Value = -0.01219512195121952
$ws.Range("D50").Value = 0.01786909370852968
$ws.Range("E50").Value = 0.001618122977346426
$ws.Range("D51").Value = 0.01742011746872979
$ws.Range("E51").Value = -0.009924385633270472
$ws.Range("D52").Value = 0.01756950848948611
$ws.Range("E52").Value = 0.005942857142857028
$ws.Range("D53").Value = 0.0167056910549516
$ws.Range("E53").Value = -0.008581937065794798
$ws.Range("D54").Value = 0.007526777436895846
$ws.Range("E54").Value = -0.008803521408563331
$ws.Range("D55").Value = 0.007545049725724912
$ws.Range("E55").Value = 0.001064509261230695
$ws.Range("E56").Value = -0.004635539119866028
